$d = $word.ActiveDocument

$pairs = @(
    @("2025-08-15 Friday", "2025-08-16 Saturday"),
    @("374÷6=", "136÷6="),
    @("332÷2=", "469÷4="),
    @("511÷5=", "361÷4="),
    @("785÷5=", "255÷3="),
    @("445÷5=", "471÷5="),
    @("513÷7=", "418÷8="),
    @("349÷9=", "939÷2="),
    @("951÷4=", "617÷7="),
    @("536÷2=", "757÷3="),
    @("754÷4=", "917÷3="),
    @("831÷3=", "895÷4="),
    @("619÷8=", "637÷2="),
    @("960÷9=", "165÷8="),
    @("906÷6=", "941÷3="),
    @("592÷4=", "574÷9="),
    @("372÷3=", "702÷2="),
    @("321÷9=", "473÷5="),
    @("238÷4=", "420÷5="),
    @("803÷7=", "890÷9="),
    @("606÷5=", "806÷6="),
    @("627÷4=", "993÷8="),
    @("522÷8=", "408÷6="),
    @("905÷3=", "947÷6="),
    @("581÷3=", "582÷5="),
    @("583÷8=", "767÷3=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
